$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3179, 3701, 3937, 4073, 4167, 4316, 4326, 4479, 4633, 4694, 4694, 4737, 4765, 5067)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("C$row").Value = $values[$i]
}
